$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 16 with the latest exam semester info
$ws.Range("A16").Value = "2024 - Vår"
$ws.Range("B16").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/hjemme-24-v.pdf)"
$ws.Range("C16").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/hjemme-24-v-solprop.html)"
$ws.Range("D16").Value = "[Materiale](tidligere-eksamensoppgaver/hjemme-24-v-ekstra.zip)"

# Update the selected cell to match the author's final selection
$ws.Range("C17").Select()
